$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# New entries (practice session dated 45916) appended as rows 349-366
# Columns: Row, B=Nom du joueur, C=Volume, D=Intensite, E=Fatigue,
#          F=Douleur, G=Localisation douleur, H=Plaisir, I=Charge (C*D)
# ---------------------------------------------------------------
$data = @(
    ,(349, 'Yoan Zouma', 70, 5, 8, 3, 'Cheville droite adduct', 0, 350)
    ,(350, 'Amir Etien', 70, 6, 8, 7, 'Coup flanc', 5, 420)
    ,(351, 'Yoann Martelat', 70, 4, 4, 5, 'Genou', 5, 280)
    ,(352, 'Jeremie Laurent', 70, 6, 5, 0, '', 5, 420)
    ,(353, 'Kamal Bafounta', 70, 5, 3, 2, 'Genou cheville', 9, 350)
    ,(354, 'Malik Boussaid', 70, 1, 0, 0, '', 10, 70)
    ,(355, 'Ilyes Boughanmi', 70, 5, 5, 4, 'Adducteur ', 10, 350)
    ,(356, 'Omar Benyounes', 70, 5, 5, 1, 'Malade', 0, 350)
    ,(357, 'Naim Ighbane', 70, 6, 6, 2, 'Cheville gauche', 4, 420)
    ,(358, 'Yanis Berrached', 70, 6, 6, 0, '', 10, 420)
    ,(359, 'Karim Belmahi', 70, 2, 7, 0, '', 10, 140)
    ,(360, 'Ilan Ihaddadene', 70, 6, 6, 0, '', 0, 420)
    ,(361, 'Emmanuel Valey', 70, 5, 2, 0, '', 8, 350)
    ,(362, 'Karahali Souaré', 70, 3, 6, 6, 'Cheville', 1, 210)
    ,(363, 'Naim Dhib', 70, 5, 5, 1, 'Genou', 0, 350)
    ,(364, 'Sofiane Belle', 70, 4, 3, 0, '', 3, 280)
    ,(365, 'Mattheo Haon', 70, 7, 5, 0, '', 2, 490)
    ,(366, 'Levy Ndoutoume', 70, 5, 5, 4, 'Cheville ischio', 2, 350)
)

$date = 45916

foreach ($d in $data) {
    $r = $d[0]
    $bName = $d[1]
    $cVal = $d[2]
    $dVal = $d[3]
    $eVal = $d[4]
    $fVal = $d[5]
    $gName = $d[6]
    $hVal = $d[7]

    # Copy the formatting of an existing row so the new cells reuse the
    # same style indices as the rest of the sheet. Row 4 has a filled
    # "Localisation douleur" cell while row 5 has an empty one.
    if ($gName -ne "") {
        $ws.Range("A4:I4").Copy() | Out-Null
    } else {
        $ws.Range("A5:I5").Copy() | Out-Null
    }
    $ws.Range("A" + $r + ":I" + $r).PasteSpecial(-4122) | Out-Null
    $excel.CutCopyMode = 0

    $ws.Range("A" + $r).Value = $date
    $ws.Range("B" + $r).Value = $bName
    $ws.Range("C" + $r).Value = $cVal
    $ws.Range("D" + $r).Value = $dVal
    $ws.Range("E" + $r).Value = $eVal
    $ws.Range("F" + $r).Value = $fVal
    if ($gName -ne "") {
        $ws.Range("G" + $r).Value = $gName
    }
    $ws.Range("H" + $r).Value = $hVal
    $ws.Range("I" + $r).Formula = "=C" + $r + "*D" + $r
}

# Widen column G ("Localisation douleur") to fit the new, longer entries
$ws.Columns.Item(7).ColumnWidth = 16

# Update the visible window: scroll near the bottom and select K360,
# matching the author's view state after the edit
$ws.Range("K360").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 335
